$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the card rows (2-5, 6-9, 10-13, 14-18, 19-25) into single
# Python-tuple-repr rows (A2:A6), then drop the now-empty trailing rows.

$ws.Range("A2").Value = "('Corrupt', ['{5}{B}', 'Sorcery', 'Corrupt deals damage to any target equal to the number of Swamps you control. You gain life equal to the damage dealt this way.'])"
$ws.Range("A3").Value = "('Gaze of Granite', ['{X}{B}{B}{G}', 'Sorcery', 'Destroy each nonland permanent with converted mana cost X or less.'])"
$ws.Range("A4").Value = "('High Tide', ['{U}', 'Instant', 'Until end of turn, whenever a player taps an Island for mana, that player adds an additional {U}.'])"
$ws.Range("A5").Value = "('Ogre Arsonist', ['{4}{R}', 'Creature — Ogre', 'When Ogre Arsonist enters the battlefield, destroy target land.', '3/3'])"
$ws.Range("A6").Value = "('Voidmage Husher', ['{3}{U}', 'Creature — Human Wizard', 'Flash (You may cast this spell any time you could cast an instant.)', 'When Voidmage Husher enters the battlefield, counter target activated ability. (Mana abilities can’t be targeted.)', 'Whenever you cast a spell, you may return Voidmage Husher to its owner’s hand.', '2/2'])"

# Remove the now-obsolete rows 7-25 entirely so the sheet's used range
# shrinks back down to A1:A6.
$ws.Range("A7:A25").ClearContents()
